# Revert "Data list npy generated": restore the tab-prefixed "ඕ" character
# in B11 (it had been re-entered without the leading tab) and restore the
# previous selection/scroll position on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = "`tඕ"

$ws.Range("C32").Select()
